# Update "Pais" sheet with refreshed per-country COVID-19 stats and re-sort
# the table by "Casos totales" (column B) descending, matching the upstream
# data refresh captured in the commit ("Update countries & provincias Spain").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowCount = 216
$data = New-Object 'object[,]' $rowCount,8

$data[0,0] = "Estados Unidos"
$data[0,1] = 5914682
$data[0,2] = 40536
$data[0,3] = 3215849
$data[0,4] = 2517739
$data[0,5] = 0
$data[0,6] = 490
$data[0,7] = 181094
$data[1,0] = "Brasil"
$data[1,1] = 3627217
$data[1,2] = 21434
$data[1,3] = 2778709
$data[1,4] = 733057
$data[1,5] = 0
$data[1,6] = 679
$data[1,7] = 115451
$data[2,0] = "India"
$data[2,1] = 3164881
$data[2,2] = 59696
$data[2,3] = 2403101
$data[2,4] = 703234
$data[2,5] = 0
$data[2,6] = 854
$data[2,7] = 58546
$data[3,0] = "Rusia"
$data[3,1] = 961493
$data[3,2] = 4744
$data[3,3] = 773095
$data[3,4] = 171950
$data[3,5] = 0
$data[3,6] = 65
$data[3,7] = 16448
$data[4,0] = "Sudafrica"
$data[4,1] = 611450
$data[4,2] = 1677
$data[4,3] = 516494
$data[4,4] = 81797
$data[4,5] = 0
$data[4,6] = 100
$data[4,7] = 13159
$data[5,0] = "Peru"
$data[5,1] = 600438
$data[5,2] = 6112
$data[5,3] = 407301
$data[5,4] = 165324
$data[5,5] = 0
$data[5,6] = 150
$data[5,7] = 27813
$data[6,0] = "Mexico"
$data[6,1] = 560164
$data[6,2] = 3948
$data[6,3] = 383872
$data[6,4] = 115812
$data[6,5] = 0
$data[6,6] = 226
$data[6,7] = 60480
$data[7,0] = "Colombia"
$data[7,1] = 551696
$data[7,2] = 10549
$data[7,3] = 384171
$data[7,4] = 149913
$data[7,5] = 0
$data[7,6] = 296
$data[7,7] = 17612
$data[8,0] = "España"
$data[8,1] = 420809
$data[8,2] = 2080
$data[8,3] = 0
$data[8,4] = 0
$data[8,5] = 0
$data[8,6] = 12
$data[8,7] = 28872
$data[9,0] = "Chile"
$data[9,1] = 399568
$data[9,2] = 1903
$data[9,3] = 372464
$data[9,4] = 16188
$data[9,5] = 0
$data[9,6] = 64
$data[9,7] = 10916
$data[10,0] = "Iran"
$data[10,1] = 361150
$data[10,2] = 2245
$data[10,3] = 311365
$data[10,4] = 29009
$data[10,5] = 0
$data[10,6] = 133
$data[10,7] = 20776
$data[11,0] = "Argentina"
$data[11,1] = 350867
$data[11,2] = 8713
$data[11,3] = 256789
$data[11,4] = 86712
$data[11,5] = 0
$data[11,6] = 381
$data[11,7] = 7366
$data[12,0] = "Reino Unido"
$data[12,1] = 326614
$data[12,2] = 853
$data[12,3] = 0
$data[12,4] = 0
$data[12,5] = 0
$data[12,6] = 4
$data[12,7] = 41433
$data[13,0] = "Arabia Saudita"
$data[13,1] = 308654
$data[13,2] = 1175
$data[13,3] = 282888
$data[13,4] = 22075
$data[13,5] = 0
$data[13,6] = 42
$data[13,7] = 3691
$data[14,0] = "Banglades"
$data[14,1] = 297083
$data[14,2] = 2485
$data[14,3] = 182875
$data[14,4] = 110225
$data[14,5] = 0
$data[14,6] = 42
$data[14,7] = 3983
$data[15,0] = "Pakistan"
$data[15,1] = 293261
$data[15,2] = 496
$data[15,3] = 276829
$data[15,4] = 10188
$data[15,5] = 0
$data[15,6] = 9
$data[15,7] = 6244
$data[16,0] = "Italia"
$data[16,1] = 260298
$data[16,2] = 953
$data[16,3] = 205662
$data[16,4] = 19195
$data[16,5] = 0
$data[16,6] = 4
$data[16,7] = 35441
$data[17,0] = "Turquia"
$data[17,1] = 259692
$data[17,2] = 1443
$data[17,3] = 237908
$data[17,4] = 15645
$data[17,5] = 0
$data[17,6] = 18
$data[17,7] = 6139
$data[18,0] = "Francia"
$data[18,1] = 244854
$data[18,2] = 1955
$data[18,3] = 85199
$data[18,4] = 129127
$data[18,5] = 0
$data[18,6] = 15
$data[18,7] = 30528
$data[19,0] = "Alemania"
$data[19,1] = 236117
$data[19,2] = 1628
$data[19,3] = 209600
$data[19,4] = 17181
$data[19,5] = 0
$data[19,6] = 4
$data[19,7] = 9336
$data[20,0] = "Irak"
$data[20,1] = 207985
$data[20,2] = 3644
$data[20,3] = 150389
$data[20,4] = 51077
$data[20,5] = 0
$data[20,6] = 91
$data[20,7] = 6519
$data[21,0] = "Filipinas"
$data[21,1] = 194252
$data[21,2] = 4686
$data[21,3] = 132042
$data[21,4] = 59200
$data[21,5] = 0
$data[21,6] = 13
$data[21,7] = 3010
$data[22,0] = "Indonesia"
$data[22,1] = 155412
$data[22,2] = 1877
$data[22,3] = 111060
$data[22,4] = 37593
$data[22,5] = 0
$data[22,6] = 79
$data[22,7] = 6759
$data[23,0] = "Canada"
$data[23,1] = 125647
$data[23,2] = 751
$data[23,3] = 111615
$data[23,4] = 4950
$data[23,5] = 0
$data[23,6] = 9
$data[23,7] = 9082
$data[24,0] = "Catar"
$data[24,1] = 117266
$data[24,2] = 258
$data[24,3] = 114099
$data[24,4] = 2973
$data[24,5] = 0
$data[24,6] = 1
$data[24,7] = 194
$data[25,0] = "Bolivia"
$data[25,1] = 109149
$data[25,2] = 722
$data[25,3] = 45396
$data[25,4] = 59244
$data[25,5] = 0
$data[25,6] = 67
$data[25,7] = 4509
$data[26,0] = "Ecuador"
$data[26,1] = 108289
$data[26,2] = 520
$data[26,3] = 94878
$data[26,4] = 7089
$data[26,5] = 0
$data[26,6] = 12
$data[26,7] = 6322
$data[27,0] = "Ucrania"
$data[27,1] = 106754
$data[27,2] = 1799
$data[27,3] = 54524
$data[27,4] = 49937
$data[27,5] = 0
$data[27,6] = 22
$data[27,7] = 2293
$data[28,0] = "Kazajistan"
$data[28,1] = 104718
$data[28,2] = 175
$data[28,3] = 92598
$data[28,4] = 10705
$data[28,5] = 0
$data[28,6] = 0
$data[28,7] = 1415
$data[29,0] = "Israel"
$data[29,1] = 104472
$data[29,2] = 1809
$data[29,3] = 81642
$data[29,4] = 21983
$data[29,5] = 0
$data[29,6] = 13
$data[29,7] = 847
$data[30,0] = "Egipto"
$data[30,1] = 97478
$data[30,2] = 138
$data[30,3] = 66817
$data[30,4] = 25381
$data[30,5] = 0
$data[30,6] = 18
$data[30,7] = 5280
$data[31,0] = "Republica Dominicana"
$data[31,1] = 91608
$data[31,2] = 447
$data[31,3] = 61558
$data[31,4] = 28477
$data[31,5] = 0
$data[31,6] = 6
$data[31,7] = 1573
$data[32,0] = "Panama"
$data[32,1] = 87485
$data[32,2] = 585
$data[32,3] = 62185
$data[32,4] = 23394
$data[32,5] = 0
$data[32,6] = 14
$data[32,7] = 1906
$data[33,0] = "Suecia"
$data[33,1] = 86721
$data[33,2] = 0
$data[33,3] = 0
$data[33,4] = 0
$data[33,5] = 0
$data[33,6] = 6
$data[33,7] = 5813
$data[34,0] = "China"
$data[34,1] = 84967
$data[34,2] = 16
$data[34,3] = 79925
$data[34,4] = 408
$data[34,5] = 0
$data[34,6] = 0
$data[34,7] = 4634
$data[35,0] = "Oman"
$data[35,1] = 84509
$data[35,2] = 740
$data[35,3] = 78912
$data[35,4] = 4960
$data[35,5] = 0
$data[35,6] = 28
$data[35,7] = 637
$data[36,0] = "Belgica"
$data[36,1] = 81936
$data[36,2] = 468
$data[36,3] = 18225
$data[36,4] = 53719
$data[36,5] = 0
$data[36,6] = 4
$data[36,7] = 9992
$data[37,0] = "Kuwait"
$data[37,1] = 80960
$data[37,2] = 432
$data[37,3] = 72925
$data[37,4] = 7517
$data[37,5] = 0
$data[37,6] = 3
$data[37,7] = 518
$data[38,0] = "Rumania"
$data[38,1] = 79330
$data[38,2] = 825
$data[38,3] = 35517
$data[38,4] = 40504
$data[38,5] = 0
$data[38,6] = 37
$data[38,7] = 3309
$data[39,0] = "Bielorrusia"
$data[39,1] = 70645
$data[39,2] = 177
$data[39,3] = 68925
$data[39,4] = 1074
$data[39,5] = 0
$data[39,6] = 4
$data[39,7] = 646
$data[40,0] = "Guatemala"
$data[40,1] = 68533
$data[40,2] = 345
$data[40,3] = 57735
$data[40,4] = 8187
$data[40,5] = 0
$data[40,6] = 17
$data[40,7] = 2611
$data[41,0] = "Emiratos Arabes Unidos"
$data[41,1] = 67282
$data[41,2] = 275
$data[41,3] = 58582
$data[41,4] = 8324
$data[41,5] = 0
$data[41,6] = 1
$data[41,7] = 376
$data[42,0] = "Paises Bajos"
$data[42,1] = 67128
$data[42,2] = 574
$data[42,3] = 0
$data[42,4] = 0
$data[42,5] = 0
$data[42,6] = 2
$data[42,7] = 6202
$data[43,0] = "Japon"
$data[43,1] = 62507
$data[43,2] = 760
$data[43,3] = 49340
$data[43,4] = 11986
$data[43,5] = 0
$data[43,6] = 5
$data[43,7] = 1181
$data[44,0] = "Polonia"
$data[44,1] = 62310
$data[44,2] = 548
$data[44,3] = 42448
$data[44,4] = 17902
$data[44,5] = 0
$data[44,6] = 5
$data[44,7] = 1960
$data[45,0] = "Singapur"
$data[45,1] = 56404
$data[45,2] = 51
$data[45,3] = 54587
$data[45,4] = 1790
$data[45,5] = 0
$data[45,6] = 0
$data[45,7] = 27
$data[46,0] = "Portugal"
$data[46,1] = 55720
$data[46,2] = 123
$data[46,3] = 40880
$data[46,4] = 13039
$data[46,5] = 0
$data[46,6] = 5
$data[46,7] = 1801
$data[47,0] = "Honduras"
$data[47,1] = 54511
$data[47,2] = 528
$data[47,3] = 8532
$data[47,4] = 44325
$data[47,5] = 0
$data[47,6] = 11
$data[47,7] = 1654
$data[48,0] = "Marruecos"
$data[48,1] = 53252
$data[48,2] = 903
$data[48,3] = 37478
$data[48,4] = 14854
$data[48,5] = 0
$data[48,6] = 32
$data[48,7] = 920
$data[49,0] = "Nigeria"
$data[49,1] = 52548
$data[49,2] = 321
$data[49,3] = 39257
$data[49,4] = 12287
$data[49,5] = 0
$data[49,6] = 2
$data[49,7] = 1004
$data[50,0] = "Barein"
$data[50,1] = 49719
$data[50,2] = 389
$data[50,3] = 46311
$data[50,4] = 3223
$data[50,5] = 0
$data[50,6] = 1
$data[50,7] = 185
$data[51,0] = "Ghana"
$data[51,1] = 43622
$data[51,2] = 117
$data[51,3] = 41695
$data[51,4] = 1664
$data[51,5] = 0
$data[51,6] = 2
$data[51,7] = 263
$data[52,0] = "Kirguistan"
$data[52,1] = 43126
$data[52,2] = 103
$data[52,3] = 36615
$data[52,4] = 5454
$data[52,5] = 0
$data[52,6] = 1
$data[52,7] = 1057
$data[53,0] = "Armenia"
$data[53,1] = 42825
$data[53,2] = 33
$data[53,3] = 36049
$data[53,4] = 5922
$data[53,5] = 0
$data[53,6] = 2
$data[53,7] = 854
$data[54,0] = "Etiopia"
$data[54,1] = 42143
$data[54,2] = 1472
$data[54,3] = 15262
$data[54,4] = 26189
$data[54,5] = 0
$data[54,6] = 14
$data[54,7] = 692
$data[55,0] = "Argelia"
$data[55,1] = 41858
$data[55,2] = 398
$data[55,3] = 29369
$data[55,4] = 11043
$data[55,5] = 0
$data[55,6] = 11
$data[55,7] = 1446
$data[56,0] = "Suiza"
$data[56,1] = 40060
$data[56,2] = 157
$data[56,3] = 34400
$data[56,4] = 3659
$data[56,5] = 0
$data[56,6] = 0
$data[56,7] = 2001
$data[57,0] = "Venezuela"
$data[57,1] = 39564
$data[57,2] = 0
$data[57,3] = 29966
$data[57,4] = 9269
$data[57,5] = 0
$data[57,6] = 0
$data[57,7] = 329
$data[58,0] = "Uzbekistan"
$data[58,1] = 39348
$data[58,2] = 402
$data[58,3] = 35551
$data[58,4] = 3517
$data[58,5] = 0
$data[58,6] = 7
$data[58,7] = 280
$data[59,0] = "Afganistan"
$data[59,1] = 38054
$data[59,2] = 55
$data[59,3] = 28360
$data[59,4] = 8305
$data[59,5] = 0
$data[59,6] = 2
$data[59,7] = 1389
$data[60,0] = "Azerbaiyan"
$data[60,1] = 35426
$data[60,2] = 152
$data[60,3] = 33104
$data[60,4] = 1803
$data[60,5] = 0
$data[60,6] = 1
$data[60,7] = 519
$data[61,0] = "Costa Rica"
$data[61,1] = 34463
$data[61,2] = 643
$data[61,3] = 12758
$data[61,4] = 21343
$data[61,5] = 0
$data[61,6] = 7
$data[61,7] = 362
$data[62,0] = "Moldavia"
$data[62,1] = 33828
$data[62,2] = 350
$data[62,3] = 23570
$data[62,4] = 9313
$data[62,5] = 0
$data[62,6] = 5
$data[62,7] = 945
$data[63,0] = "Nepal"
$data[63,1] = 32678
$data[63,2] = 743
$data[63,3] = 18806
$data[63,4] = 13715
$data[63,5] = 0
$data[63,6] = 8
$data[63,7] = 157
$data[64,0] = "Kenia"
$data[64,1] = 32557
$data[64,2] = 193
$data[64,3] = 18895
$data[64,4] = 13108
$data[64,5] = 0
$data[64,6] = 6
$data[64,7] = 554
$data[65,0] = "Serbia"
$data[65,1] = 30714
$data[65,2] = 57
$data[65,3] = 29028
$data[65,4] = 985
$data[65,5] = 0
$data[65,6] = 3
$data[65,7] = 701
$data[66,0] = "Irlanda"
$data[66,1] = 28116
$data[66,2] = 147
$data[66,3] = 23364
$data[66,4] = 2975
$data[66,5] = 0
$data[66,6] = 0
$data[66,7] = 1777
$data[67,0] = "Austria"
$data[67,1] = 25495
$data[67,2] = 242
$data[67,3] = 21657
$data[67,4] = 3105
$data[67,5] = 0
$data[67,6] = 1
$data[67,7] = 733
$data[68,0] = "Australia"
$data[68,1] = 24916
$data[68,2] = 104
$data[68,3] = 19603
$data[68,4] = 4796
$data[68,5] = 0
$data[68,6] = 15
$data[68,7] = 517
$data[69,0] = "El Salvador"
$data[69,1] = 24811
$data[69,2] = 189
$data[69,3] = 12492
$data[69,4] = 11650
$data[69,5] = 0
$data[69,6] = 8
$data[69,7] = 669
$data[70,0] = "Chequia"
$data[70,1] = 22181
$data[70,2] = 258
$data[70,3] = 16376
$data[70,4] = 5390
$data[70,5] = 0
$data[70,6] = 3
$data[70,7] = 415
$data[71,0] = "Estado de Palestina"
$data[71,1] = 19213
$data[71,2] = 411
$data[71,3] = 11870
$data[71,4] = 7210
$data[71,5] = 0
$data[71,6] = 5
$data[71,7] = 133
$data[72,0] = "Camerun"
$data[72,1] = 18762
$data[72,2] = 0
$data[72,3] = 16540
$data[72,4] = 1814
$data[72,5] = 0
$data[72,6] = 0
$data[72,7] = 408
$data[73,0] = "Bosnia y Herzegovina"
$data[73,1] = 18029
$data[73,2] = 314
$data[73,3] = 11861
$data[73,4] = 5621
$data[73,5] = 0
$data[73,6] = 15
$data[73,7] = 547
$data[74,0] = "Corea del Sur"
$data[74,1] = 17665
$data[74,2] = 266
$data[74,3] = 14219
$data[74,4] = 3137
$data[74,5] = 0
$data[74,6] = 0
$data[74,7] = 309
$data[75,0] = "Costa de Marfil"
$data[75,1] = 17506
$data[75,2] = 35
$data[75,3] = 15633
$data[75,4] = 1759
$data[75,5] = 0
$data[75,6] = 1
$data[75,7] = 114
$data[76,0] = "Dinamarca"
$data[76,1] = 16397
$data[76,2] = 80
$data[76,3] = 14310
$data[76,4] = 1464
$data[76,5] = 0
$data[76,6] = 1
$data[76,7] = 623
$data[77,0] = "Bulgaria"
$data[77,1] = 15386
$data[77,2] = 99
$data[77,3] = 10497
$data[77,4] = 4326
$data[77,5] = 0
$data[77,6] = 18
$data[77,7] = 563
$data[78,0] = "Madagascar"
$data[78,1] = 14402
$data[78,2] = 75
$data[78,3] = 13436
$data[78,4] = 788
$data[78,5] = 0
$data[78,6] = 0
$data[78,7] = 178
$data[79,0] = "Republica de Macedonia"
$data[79,1] = 13673
$data[79,2] = 78
$data[79,3] = 10150
$data[79,4] = 2955
$data[79,5] = 0
$data[79,6] = 4
$data[79,7] = 568
$data[80,0] = "Paraguay"
$data[80,1] = 13233
$data[80,2] = 0
$data[80,3] = 7417
$data[80,4] = 5611
$data[80,5] = 0
$data[80,6] = 0
$data[80,7] = 205
$data[81,0] = "Libano"
$data[81,1] = 13155
$data[81,2] = 457
$data[81,3] = 3704
$data[81,4] = 9325
$data[81,5] = 0
$data[81,6] = 3
$data[81,7] = 126
$data[82,0] = "Senegal"
$data[82,1] = 13013
$data[82,2] = 64
$data[82,3] = 8595
$data[82,4] = 4146
$data[82,5] = 0
$data[82,6] = 3
$data[82,7] = 272
$data[83,0] = "Sudan"
$data[83,1] = 12903
$data[83,2] = 67
$data[83,3] = 6538
$data[83,4] = 5547
$data[83,5] = 0
$data[83,6] = 3
$data[83,7] = 818
$data[84,0] = "Zambia"
$data[84,1] = 11148
$data[84,2] = 66
$data[84,3] = 10208
$data[84,4] = 660
$data[84,5] = 0
$data[84,6] = 0
$data[84,7] = 280
$data[85,0] = "Libia"
$data[85,1] = 11009
$data[85,2] = 572
$data[85,3] = 1096
$data[85,4] = 9714
$data[85,5] = 0
$data[85,6] = 11
$data[85,7] = 199
$data[86,0] = "Noruega"
$data[86,1] = 10395
$data[86,2] = 72
$data[86,3] = 9150
$data[86,4] = 981
$data[86,5] = 0
$data[86,6] = 0
$data[86,7] = 264
$data[87,0] = "Consejo Danes para los Refugiados"
$data[87,1] = 9842
$data[87,2] = 12
$data[87,3] = 8953
$data[87,4] = 638
$data[87,5] = 0
$data[87,6] = 0
$data[87,7] = 251
$data[88,0] = "Malasia"
$data[88,1] = 9274
$data[88,2] = 7
$data[88,3] = 8965
$data[88,4] = 184
$data[88,5] = 0
$data[88,6] = 0
$data[88,7] = 125
$data[89,0] = "Guinea"
$data[89,1] = 9076
$data[89,2] = 109
$data[89,3] = 7928
$data[89,4] = 1093
$data[89,5] = 0
$data[89,6] = 2
$data[89,7] = 55
$data[90,0] = "Guayana Francesa"
$data[90,1] = 8875
$data[90,2] = 78
$data[90,3] = 8363
$data[90,4] = 456
$data[90,5] = 0
$data[90,6] = 1
$data[90,7] = 56
$data[91,0] = "Grecia"
$data[91,1] = 8819
$data[91,2] = 155
$data[91,3] = 3804
$data[91,4] = 4773
$data[91,5] = 0
$data[91,6] = 0
$data[91,7] = 242
$data[92,0] = "Albania"
$data[92,1] = 8605
$data[92,2] = 178
$data[92,3] = 4413
$data[92,4] = 3938
$data[92,5] = 0
$data[92,6] = 4
$data[92,7] = 254
$data[93,0] = "Gabon"
$data[93,1] = 8409
$data[93,2] = 21
$data[93,3] = 6959
$data[93,4] = 1397
$data[93,5] = 0
$data[93,6] = 0
$data[93,7] = 53
$data[94,0] = "Tayikistan"
$data[94,1] = 8346
$data[94,2] = 35
$data[94,3] = 7142
$data[94,4] = 1137
$data[94,5] = 0
$data[94,6] = 1
$data[94,7] = 67
$data[95,0] = "Croacia"
$data[95,1] = 8311
$data[95,2] = 136
$data[95,3] = 5926
$data[95,4] = 2212
$data[95,5] = 0
$data[95,6] = 2
$data[95,7] = 173
$data[96,0] = "Haiti"
$data[96,1] = 8110
$data[96,2] = 28
$data[96,3] = 5624
$data[96,4] = 2290
$data[96,5] = 0
$data[96,6] = 0
$data[96,7] = 196
$data[97,0] = "Finlandia"
$data[97,1] = 7938
$data[97,2] = 18
$data[97,3] = 7100
$data[97,4] = 503
$data[97,5] = 0
$data[97,6] = 1
$data[97,7] = 335
$data[98,0] = "Luxemburgo"
$data[98,1] = 7794
$data[98,2] = 19
$data[98,3] = 7106
$data[98,4] = 564
$data[98,5] = 0
$data[98,6] = 0
$data[98,7] = 124
$data[99,0] = "Mauritania"
$data[99,1] = 6928
$data[99,2] = 23
$data[99,3] = 6282
$data[99,4] = 488
$data[99,5] = 0
$data[99,6] = 0
$data[99,7] = 158
$data[100,0] = "Maldivas"
$data[100,1] = 6912
$data[100,2] = 133
$data[100,3] = 4297
$data[100,4] = 2588
$data[100,5] = 0
$data[100,6] = 1
$data[100,7] = 27
$data[101,0] = "Zimbabue"
$data[101,1] = 6070
$data[101,2] = 140
$data[101,3] = 4950
$data[101,4] = 965
$data[101,5] = 0
$data[101,6] = 0
$data[101,7] = 155
$data[102,0] = "Namibia"
$data[102,1] = 6030
$data[102,2] = 176
$data[102,3] = 2563
$data[102,4] = 3411
$data[102,5] = 0
$data[102,6] = 4
$data[102,7] = 56
$data[103,0] = "Malaui"
$data[103,1] = 5419
$data[103,2] = 5
$data[103,3] = 3059
$data[103,4] = 2191
$data[103,5] = 0
$data[103,6] = 1
$data[103,7] = 169
$data[104,0] = "Republica de Yibuti"
$data[104,1] = 5383
$data[104,2] = 1
$data[104,3] = 5273
$data[104,4] = 50
$data[104,5] = 0
$data[104,6] = 0
$data[104,7] = 60
$data[105,0] = "Hungria"
$data[105,1] = 5191
$data[105,2] = 36
$data[105,3] = 3695
$data[105,4] = 883
$data[105,5] = 0
$data[105,6] = 0
$data[105,7] = 613
$data[106,0] = "Guinea Ecuatorial"
$data[106,1] = 4926
$data[106,2] = 0
$data[106,3] = 3795
$data[106,4] = 1048
$data[106,5] = 0
$data[106,6] = 0
$data[106,7] = 83
$data[107,0] = "Hong Kong"
$data[107,1] = 4692
$data[107,2] = 9
$data[107,3] = 4052
$data[107,4] = 563
$data[107,5] = 0
$data[107,6] = 0
$data[107,7] = 77
$data[108,0] = "Republica de Africa Central"
$data[108,1] = 4679
$data[108,2] = 0
$data[108,3] = 1755
$data[108,4] = 2863
$data[108,5] = 0
$data[108,6] = 0
$data[108,7] = 61
$data[109,0] = "Montenegro"
$data[109,1] = 4378
$data[109,2] = 35
$data[109,3] = 3420
$data[109,4] = 874
$data[109,5] = 0
$data[109,6] = 0
$data[109,7] = 84
$data[110,0] = "Nicaragua"
$data[110,1] = 4311
$data[110,2] = 0
$data[110,3] = 2913
$data[110,4] = 1265
$data[110,5] = 0
$data[110,6] = 0
$data[110,7] = 133
$data[111,0] = "Suazilandia"
$data[111,1] = 4304
$data[111,2] = 79
$data[111,3] = 2936
$data[111,4] = 1283
$data[111,5] = 0
$data[111,6] = 0
$data[111,7] = 85
$data[112,0] = "Congo"
$data[112,1] = 3979
$data[112,2] = 129
$data[112,3] = 1742
$data[112,4] = 2159
$data[112,5] = 0
$data[112,6] = 1
$data[112,7] = 78
$data[113,0] = "Cuba"
$data[113,1] = 3717
$data[113,2] = 35
$data[113,3] = 3079
$data[113,4] = 547
$data[113,5] = 0
$data[113,6] = 0
$data[113,7] = 91
$data[114,0] = "Surinam"
$data[114,1] = 3632
$data[114,2] = 25
$data[114,3] = 2758
$data[114,4] = 814
$data[114,5] = 0
$data[114,6] = 2
$data[114,7] = 60
$data[115,0] = "Cabo Verde"
$data[115,1] = 3532
$data[115,2] = 23
$data[115,3] = 2599
$data[115,4] = 896
$data[115,5] = 0
$data[115,6] = 0
$data[115,7] = 37
$data[116,0] = "Mozambique"
$data[116,1] = 3440
$data[116,2] = 45
$data[116,3] = 1661
$data[116,4] = 1758
$data[116,5] = 0
$data[116,6] = 1
$data[116,7] = 21
$data[117,0] = "Eslovaquia"
$data[117,1] = 3424
$data[117,2] = 68
$data[117,3] = 2153
$data[117,4] = 1238
$data[117,5] = 0
$data[117,6] = 0
$data[117,7] = 33
$data[118,0] = "Tailandia"
$data[118,1] = 3397
$data[118,2] = 2
$data[118,3] = 3222
$data[118,4] = 117
$data[118,5] = 0
$data[118,6] = 0
$data[118,7] = 58
$data[119,0] = "Ruanda"
$data[119,1] = 3306
$data[119,2] = 217
$data[119,3] = 1785
$data[119,4] = 1507
$data[119,5] = 0
$data[119,6] = 2
$data[119,7] = 14
$data[120,0] = "Somalia"
$data[120,1] = 3269
$data[120,2] = 0
$data[120,3] = 2443
$data[120,4] = 733
$data[120,5] = 0
$data[120,6] = 0
$data[120,7] = 93
$data[121,0] = "Mayotte"
$data[121,1] = 3237
$data[121,2] = 0
$data[121,3] = 2964
$data[121,4] = 234
$data[121,5] = 0
$data[121,6] = 0
$data[121,7] = 39
$data[122,0] = "Sri Lanka"
$data[122,1] = 2959
$data[122,2] = 6
$data[122,3] = 2811
$data[122,4] = 136
$data[122,5] = 0
$data[122,6] = 0
$data[122,7] = 12
$data[123,0] = "Tunez"
$data[123,1] = 2893
$data[123,2] = 75
$data[123,3] = 1454
$data[123,4] = 1368
$data[123,5] = 0
$data[123,6] = 0
$data[123,7] = 71
$data[124,0] = "Mali"
$data[124,1] = 2708
$data[124,2] = 3
$data[124,3] = 2025
$data[124,4] = 558
$data[124,5] = 0
$data[124,6] = 0
$data[124,7] = 125
$data[125,0] = "Lituania"
$data[125,1] = 2673
$data[125,2] = 38
$data[125,3] = 1766
$data[125,4] = 822
$data[125,5] = 0
$data[125,6] = 1
$data[125,7] = 85
$data[126,0] = "Eslovenia"
$data[126,1] = 2665
$data[126,2] = 14
$data[126,3] = 2122
$data[126,4] = 410
$data[126,5] = 0
$data[126,6] = 2
$data[126,7] = 133
$data[127,0] = "Gambia"
$data[127,1] = 2585
$data[127,2] = 148
$data[127,3] = 490
$data[127,4] = 2008
$data[127,5] = 0
$data[127,6] = 3
$data[127,7] = 87
$data[128,0] = "Sudan del Sur"
$data[128,1] = 2504
$data[128,2] = 5
$data[128,3] = 1290
$data[128,4] = 1167
$data[128,5] = 0
$data[128,6] = 0
$data[128,7] = 47
$data[129,0] = "Uganda"
$data[129,1] = 2362
$data[129,2] = 99
$data[129,3] = 1248
$data[129,4] = 1092
$data[129,5] = 0
$data[129,6] = 2
$data[129,7] = 22
$data[130,0] = "Siria"
$data[130,1] = 2293
$data[130,2] = 76
$data[130,3] = 519
$data[130,4] = 1682
$data[130,5] = 0
$data[130,6] = 3
$data[130,7] = 92
$data[131,0] = "Estonia"
$data[131,1] = 2275
$data[131,2] = 3
$data[131,3] = 2025
$data[131,4] = 186
$data[131,5] = 0
$data[131,6] = 1
$data[131,7] = 64
$data[132,0] = "Angola"
$data[132,1] = 2222
$data[132,2] = 51
$data[132,3] = 877
$data[132,4] = 1245
$data[132,5] = 0
$data[132,6] = 4
$data[132,7] = 100
$data[133,0] = "Guinea-Bisau"
$data[133,1] = 2149
$data[133,2] = 0
$data[133,3] = 1015
$data[133,4] = 1101
$data[133,5] = 0
$data[133,6] = 0
$data[133,7] = 33
$data[134,0] = "Benin"
$data[134,1] = 2115
$data[134,2] = 0
$data[134,3] = 1705
$data[134,4] = 371
$data[134,5] = 0
$data[134,6] = 0
$data[134,7] = 39
$data[135,0] = "Islandia"
$data[135,1] = 2073
$data[135,2] = 9
$data[135,3] = 1946
$data[135,4] = 117
$data[135,5] = 0
$data[135,6] = 0
$data[135,7] = 10
$data[136,0] = "Sierra Leona"
$data[136,1] = 1997
$data[136,2] = 5
$data[136,3] = 1557
$data[136,4] = 371
$data[136,5] = 0
$data[136,6] = 0
$data[136,7] = 69
$data[137,0] = "Yemen"
$data[137,1] = 1916
$data[137,2] = 5
$data[137,3] = 1090
$data[137,4] = 271
$data[137,5] = 0
$data[137,6] = 2
$data[137,7] = 555
$data[138,0] = "Bahamas"
$data[138,1] = 1765
$data[138,2] = 0
$data[138,3] = 227
$data[138,4] = 1509
$data[138,5] = 0
$data[138,6] = 0
$data[138,7] = 29
$data[139,0] = "Nueva Zelanda"
$data[139,1] = 1683
$data[139,2] = 9
$data[139,3] = 1538
$data[139,4] = 123
$data[139,5] = 0
$data[139,6] = 0
$data[139,7] = 22
$data[140,0] = "Malta"
$data[140,1] = 1667
$data[140,2] = 55
$data[140,3] = 977
$data[140,4] = 680
$data[140,5] = 0
$data[140,6] = 0
$data[140,7] = 10
$data[141,0] = "Jordania"
$data[141,1] = 1639
$data[141,2] = 30
$data[141,3] = 1335
$data[141,4] = 290
$data[141,5] = 0
$data[141,6] = 2
$data[141,7] = 14
$data[142,0] = "Aruba"
$data[142,1] = 1628
$data[142,2] = 60
$data[142,3] = 461
$data[142,4] = 1160
$data[142,5] = 0
$data[142,6] = 0
$data[142,7] = 7
$data[143,0] = "Botsuana"
$data[143,1] = 1562
$data[143,2] = 254
$data[143,3] = 199
$data[143,4] = 1360
$data[143,5] = 0
$data[143,6] = 0
$data[143,7] = 3
$data[144,0] = "Uruguay"
$data[144,1] = 1533
$data[144,2] = 6
$data[144,3] = 1295
$data[144,4] = 196
$data[144,5] = 0
$data[144,6] = 0
$data[144,7] = 42
$data[145,0] = "Jamaica"
$data[145,1] = 1529
$data[145,2] = 116
$data[145,3] = 819
$data[145,4] = 694
$data[145,5] = 0
$data[145,6] = 0
$data[145,7] = 16
$data[146,0] = "Republica de Chipre"
$data[146,1] = 1451
$data[146,2] = 30
$data[146,3] = 878
$data[146,4] = 553
$data[146,5] = 0
$data[146,6] = 0
$data[146,7] = 20
$data[147,0] = "Georgia"
$data[147,1] = 1421
$data[147,2] = 10
$data[147,3] = 1137
$data[147,4] = 266
$data[147,5] = 0
$data[147,6] = 1
$data[147,7] = 18
$data[148,0] = "Burkina Faso"
$data[148,1] = 1338
$data[148,2] = 18
$data[148,3] = 1050
$data[148,4] = 233
$data[148,5] = 0
$data[148,6] = 0
$data[148,7] = 55
$data[149,0] = "Letonia"
$data[149,1] = 1337
$data[149,2] = 0
$data[149,3] = 1093
$data[149,4] = 211
$data[149,5] = 0
$data[149,6] = 0
$data[149,7] = 33
$data[150,0] = "Togo"
$data[150,1] = 1295
$data[150,2] = 18
$data[150,3] = 914
$data[150,4] = 354
$data[150,5] = 0
$data[150,6] = 0
$data[150,7] = 27
$data[151,0] = "Liberia"
$data[151,1] = 1290
$data[151,2] = 4
$data[151,3] = 819
$data[151,4] = 389
$data[151,5] = 0
$data[151,6] = 0
$data[151,7] = 82
$data[152,0] = "Reunion"
$data[152,1] = 1244
$data[152,2] = 35
$data[152,3] = 692
$data[152,4] = 546
$data[152,5] = 0
$data[152,6] = 0
$data[152,7] = 6
$data[153,0] = "Niger"
$data[153,1] = 1172
$data[153,2] = 0
$data[153,3] = 1084
$data[153,4] = 19
$data[153,5] = 0
$data[153,6] = 0
$data[153,7] = 69
$data[154,0] = "Trinidad yTobago"
$data[154,1] = 1099
$data[154,2] = 92
$data[154,3] = 165
$data[154,4] = 919
$data[154,5] = 0
$data[154,6] = 1
$data[154,7] = 15
$data[155,0] = "Principado de Andorra"
$data[155,1] = 1060
$data[155,2] = 15
$data[155,3] = 877
$data[155,4] = 130
$data[155,5] = 0
$data[155,6] = 0
$data[155,7] = 53
$data[156,0] = "Guyana"
$data[156,1] = 1029
$data[156,2] = 74
$data[156,3] = 510
$data[156,4] = 488
$data[156,5] = 0
$data[156,6] = 0
$data[156,7] = 31
$data[157,0] = "Vietnam"
$data[157,1] = 1022
$data[157,2] = 6
$data[157,3] = 587
$data[157,4] = 408
$data[157,5] = 0
$data[157,6] = 0
$data[157,7] = 27
$data[158,0] = "Lesoto"
$data[158,1] = 1015
$data[158,2] = 0
$data[158,3] = 472
$data[158,4] = 513
$data[158,5] = 0
$data[158,6] = 0
$data[158,7] = 30
$data[159,0] = "Republica del Chad"
$data[159,1] = 987
$data[159,2] = 1
$data[159,3] = 870
$data[159,4] = 41
$data[159,5] = 0
$data[159,6] = 0
$data[159,7] = 76
$data[160,0] = "Santo Tome y Principe"
$data[160,1] = 892
$data[160,2] = 0
$data[160,3] = 831
$data[160,4] = 46
$data[160,5] = 0
$data[160,6] = 0
$data[160,7] = 15
$data[161,0] = "Guadalupe"
$data[161,1] = 771
$data[161,2] = 0
$data[161,3] = 289
$data[161,4] = 467
$data[161,5] = 0
$data[161,6] = 0
$data[161,7] = 15
$data[162,0] = "Belice"
$data[162,1] = 713
$data[162,2] = 27
$data[162,3] = 45
$data[162,4] = 658
$data[162,5] = 0
$data[162,6] = 4
$data[162,7] = 10
$data[163,0] = "Crucero"
$data[163,1] = 712
$data[163,2] = 0
$data[163,3] = 651
$data[163,4] = 48
$data[163,5] = 0
$data[163,6] = 0
$data[163,7] = 13
$data[164,0] = "San Marino"
$data[164,1] = 704
$data[164,2] = 0
$data[164,3] = 657
$data[164,4] = 5
$data[164,5] = 0
$data[164,6] = 0
$data[164,7] = 42
$data[165,0] = "Tanzania"
$data[165,1] = 509
$data[165,2] = 0
$data[165,3] = 183
$data[165,4] = 305
$data[165,5] = 0
$data[165,6] = 0
$data[165,7] = 21
$data[166,0] = "Taiwan"
$data[166,1] = 487
$data[166,2] = 0
$data[166,3] = 457
$data[166,4] = 23
$data[166,5] = 0
$data[166,6] = 0
$data[166,7] = 7
$data[167,0] = "Birmania"
$data[167,1] = 474
$data[167,2] = 24
$data[167,3] = 341
$data[167,4] = 127
$data[167,5] = 0
$data[167,6] = 0
$data[167,7] = 6
$data[168,0] = "Martinica"
$data[168,1] = 464
$data[168,2] = 0
$data[168,3] = 98
$data[168,4] = 350
$data[168,5] = 0
$data[168,6] = 0
$data[168,7] = 16
$data[169,0] = "Burundi"
$data[169,1] = 430
$data[169,2] = 0
$data[169,3] = 336
$data[169,4] = 93
$data[169,5] = 0
$data[169,6] = 0
$data[169,7] = 1
$data[170,0] = "Comoras"
$data[170,1] = 417
$data[170,2] = 0
$data[170,3] = 396
$data[170,4] = 14
$data[170,5] = 0
$data[170,6] = 0
$data[170,7] = 7
$data[171,0] = "Islas Feroe"
$data[171,1] = 410
$data[171,2] = 26
$data[171,3] = 337
$data[171,4] = 73
$data[171,5] = 0
$data[171,6] = 0
$data[171,7] = 0
$data[172,0] = "Papua Nueva Guinea"
$data[172,1] = 401
$data[172,2] = 0
$data[172,3] = 232
$data[172,4] = 165
$data[172,5] = 0
$data[172,6] = 0
$data[172,7] = 4
$data[173,0] = "San Martin (Parte Holandesa)"
$data[173,1] = 396
$data[173,2] = 0
$data[173,3] = 147
$data[173,4] = 232
$data[173,5] = 0
$data[173,6] = 0
$data[173,7] = 17
$data[174,0] = "Islas Turcas y Caicos"
$data[174,1] = 383
$data[174,2] = 0
$data[174,3] = 102
$data[174,4] = 279
$data[174,5] = 0
$data[174,6] = 0
$data[174,7] = 2
$data[175,0] = "Mauricio"
$data[175,1] = 347
$data[175,2] = 1
$data[175,3] = 335
$data[175,4] = 2
$data[175,5] = 0
$data[175,6] = 0
$data[175,7] = 10
$data[176,0] = "Isla de Man"
$data[176,1] = 336
$data[176,2] = 0
$data[176,3] = 312
$data[176,4] = 0
$data[176,5] = 0
$data[176,6] = 0
$data[176,7] = 24
$data[177,0] = "Eritrea"
$data[177,1] = 306
$data[177,2] = 0
$data[177,3] = 274
$data[177,4] = 32
$data[177,5] = 0
$data[177,6] = 0
$data[177,7] = 0
$data[178,0] = "Polinesia Francesa"
$data[178,1] = 298
$data[178,2] = 0
$data[178,3] = 148
$data[178,4] = 150
$data[178,5] = 0
$data[178,6] = 0
$data[178,7] = 0
$data[179,0] = "Mongolia"
$data[179,1] = 298
$data[179,2] = 0
$data[179,3] = 289
$data[179,4] = 9
$data[179,5] = 0
$data[179,6] = 0
$data[179,7] = 0
$data[180,0] = "Camboya"
$data[180,1] = 273
$data[180,2] = 0
$data[180,3] = 263
$data[180,4] = 10
$data[180,5] = 0
$data[180,6] = 0
$data[180,7] = 0
$data[181,0] = "Gibraltar"
$data[181,1] = 248
$data[181,2] = 2
$data[181,3] = 203
$data[181,4] = 45
$data[181,5] = 0
$data[181,6] = 0
$data[181,7] = 0
$data[182,0] = "Islas Caimanes"
$data[182,1] = 205
$data[182,2] = 0
$data[182,3] = 202
$data[182,4] = 2
$data[182,5] = 0
$data[182,6] = 0
$data[182,7] = 1
$data[183,0] = "San Martin (Parte Francesa)"
$data[183,1] = 176
$data[183,2] = 0
$data[183,3] = 52
$data[183,4] = 119
$data[183,5] = 0
$data[183,6] = 0
$data[183,7] = 5
$data[184,0] = "Bermudas"
$data[184,1] = 167
$data[184,2] = 0
$data[184,3] = 149
$data[184,4] = 9
$data[184,5] = 0
$data[184,6] = 0
$data[184,7] = 9
$data[185,0] = "Barbados"
$data[185,1] = 161
$data[185,2] = 0
$data[185,3] = 132
$data[185,4] = 22
$data[185,5] = 0
$data[185,6] = 0
$data[185,7] = 7
$data[186,0] = "Butan"
$data[186,1] = 155
$data[186,2] = 0
$data[186,3] = 115
$data[186,4] = 40
$data[186,5] = 0
$data[186,6] = 0
$data[186,7] = 0
$data[187,0] = "Brunei"
$data[187,1] = 143
$data[187,2] = 0
$data[187,3] = 139
$data[187,4] = 1
$data[187,5] = 0
$data[187,6] = 0
$data[187,7] = 3
$data[188,0] = "Seychelles"
$data[188,1] = 132
$data[188,2] = 0
$data[188,3] = 126
$data[188,4] = 6
$data[188,5] = 0
$data[188,6] = 0
$data[188,7] = 0
$data[189,0] = "Monaco"
$data[189,1] = 112
$data[189,2] = 0
$data[189,3] = 83
$data[189,4] = 28
$data[189,5] = 0
$data[189,6] = 0
$data[189,7] = 1
$data[190,0] = "Liechtenstein"
$data[190,1] = 100
$data[190,2] = 1
$data[190,3] = 90
$data[190,4] = 9
$data[190,5] = 0
$data[190,6] = 0
$data[190,7] = 1
$data[191,0] = "Antigua y Barbuda"
$data[191,1] = 94
$data[191,2] = 0
$data[191,3] = 89
$data[191,4] = 2
$data[191,5] = 0
$data[191,6] = 0
$data[191,7] = 3
$data[192,0] = "San Vicente y las Granadinas"
$data[192,1] = 58
$data[192,2] = 0
$data[192,3] = 56
$data[192,4] = 2
$data[192,5] = 0
$data[192,6] = 0
$data[192,7] = 0
$data[193,0] = "Macao"
$data[193,1] = 46
$data[193,2] = 0
$data[193,3] = 46
$data[193,4] = 0
$data[193,5] = 0
$data[193,6] = 0
$data[193,7] = 0
$data[194,0] = "Curazao"
$data[194,1] = 43
$data[194,2] = 4
$data[194,3] = 34
$data[194,4] = 8
$data[194,5] = 0
$data[194,6] = 0
$data[194,7] = 1
$data[195,0] = "Puerto Rico"
$data[195,1] = 39
$data[195,2] = 0
$data[195,3] = 1
$data[195,4] = 36
$data[195,5] = 0
$data[195,6] = 0
$data[195,7] = 2
$data[196,0] = "Guam"
$data[196,1] = 32
$data[196,2] = 0
$data[196,3] = 0
$data[196,4] = 31
$data[196,5] = 0
$data[196,6] = 0
$data[196,7] = 1
$data[197,0] = "Fiyi"
$data[197,1] = 28
$data[197,2] = 0
$data[197,3] = 23
$data[197,4] = 4
$data[197,5] = 0
$data[197,6] = 0
$data[197,7] = 1
$data[198,0] = "Timor Oriental"
$data[198,1] = 26
$data[198,2] = 0
$data[198,3] = 25
$data[198,4] = 1
$data[198,5] = 0
$data[198,6] = 0
$data[198,7] = 0
$data[199,0] = "Santa Lucia"
$data[199,1] = 26
$data[199,2] = 0
$data[199,3] = 25
$data[199,4] = 1
$data[199,5] = 0
$data[199,6] = 0
$data[199,7] = 0
$data[200,0] = "Granada"
$data[200,1] = 24
$data[200,2] = 0
$data[200,3] = 24
$data[200,4] = 0
$data[200,5] = 0
$data[200,6] = 0
$data[200,7] = 0
$data[201,0] = "Nueva Caledonia"
$data[201,1] = 23
$data[201,2] = 0
$data[201,3] = 23
$data[201,4] = 0
$data[201,5] = 0
$data[201,6] = 0
$data[201,7] = 0
$data[202,0] = "Laos"
$data[202,1] = 22
$data[202,2] = 0
$data[202,3] = 20
$data[202,4] = 2
$data[202,5] = 0
$data[202,6] = 0
$data[202,7] = 0
$data[203,0] = "Islas Virgenes Britanicas"
$data[203,1] = 21
$data[203,2] = 0
$data[203,3] = 8
$data[203,4] = 12
$data[203,5] = 0
$data[203,6] = 0
$data[203,7] = 1
$data[204,0] = "Dominica"
$data[204,1] = 20
$data[204,2] = 1
$data[204,3] = 18
$data[204,4] = 2
$data[204,5] = 0
$data[204,6] = 0
$data[204,7] = 0
$data[205,0] = "Islas Virgenes de los Estados Unidos"
$data[205,1] = 17
$data[205,2] = 0
$data[205,3] = 0
$data[205,4] = 17
$data[205,5] = 0
$data[205,6] = 0
$data[205,7] = 0
$data[206,0] = "San Cristobal y Nieves"
$data[206,1] = 17
$data[206,2] = 0
$data[206,3] = 17
$data[206,4] = 0
$data[206,5] = 0
$data[206,6] = 0
$data[206,7] = 0
$data[207,0] = "San Bartolome"
$data[207,1] = 16
$data[207,2] = 0
$data[207,3] = 9
$data[207,4] = 7
$data[207,5] = 0
$data[207,6] = 0
$data[207,7] = 0
$data[208,0] = "Groenlandia"
$data[208,1] = 14
$data[208,2] = 0
$data[208,3] = 14
$data[208,4] = 0
$data[208,5] = 0
$data[208,6] = 0
$data[208,7] = 0
$data[209,0] = "Bonaire, San Eustaquio y Saba"
$data[209,1] = 13
$data[209,2] = 0
$data[209,3] = 7
$data[209,4] = 6
$data[209,5] = 0
$data[209,6] = 0
$data[209,7] = 0
$data[210,0] = "Islas Malvinas"
$data[210,1] = 13
$data[210,2] = 0
$data[210,3] = 13
$data[210,4] = 0
$data[210,5] = 0
$data[210,6] = 0
$data[210,7] = 0
$data[211,0] = "Montserrat"
$data[211,1] = 13
$data[211,2] = 0
$data[211,3] = 12
$data[211,4] = 0
$data[211,5] = 0
$data[211,6] = 0
$data[211,7] = 1
$data[212,0] = "Santa Sede"
$data[212,1] = 12
$data[212,2] = 0
$data[212,3] = 12
$data[212,4] = 0
$data[212,5] = 0
$data[212,6] = 0
$data[212,7] = 0
$data[213,0] = "Sahara Occidental"
$data[213,1] = 10
$data[213,2] = 0
$data[213,3] = 8
$data[213,4] = 1
$data[213,5] = 0
$data[213,6] = 0
$data[213,7] = 1
$data[214,0] = "San Pedro y Miquelon"
$data[214,1] = 5
$data[214,2] = 0
$data[214,3] = 1
$data[214,4] = 4
$data[214,5] = 0
$data[214,6] = 0
$data[214,7] = 0
$data[215,0] = "Anguila"
$data[215,1] = 3
$data[215,2] = 0
$data[215,3] = 3
$data[215,4] = 0
$data[215,5] = 0
$data[215,6] = 0
$data[215,7] = 0

# Countries occupy A4:H219 (row 1 = updated-at banner, row 3 = header).
$ws.Range("A4:H219").Value = $data

# Refresh the "last updated" banner in A1.
$ws.Range("A1").Value = "Datos actualizados a 25 de Agosto de 2020 a las 02:18"
